$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1, matching the style of the existing header (E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# Boolean (MAD outlier) values for rows 2-20 across columns F (KNN), G (SVM), H (RF)
$values = @{
    2  = @($false, $false, $false)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($true,  $true,  $true)
    6  = @($true,  $false, $true)
    7  = @($false, $false, $false)
    8  = @($false, $false, $false)
    9  = @($false, $false, $false)
    10 = @($true,  $true,  $true)
    11 = @($false, $false, $false)
    12 = @($false, $false, $false)
    13 = @($true,  $false, $false)
    14 = @($false, $false, $false)
    15 = @($false, $false, $false)
    16 = @($false, $false, $true)
    17 = @($false, $false, $false)
    18 = @($false, $false, $false)
    19 = @($false, $false, $false)
    20 = @($false, $false, $false)
}

foreach ($row in $values.Keys | Sort-Object) {
    $vals = $values[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}

$wb.Save()
